$d = $word.ActiveDocument

# 1) Replace the ID placeholder text, absorbing the trailing space run into
#    the replacement so only a single run with the new text remains.
[void]$d.Content.Find.Execute("**ID__AFFARS_5309_topic_13__ID** ", $true, $false, $false, $false, $false, $true, 1, $false, "**ID__AFFARS_5309_405_1__ID**", 2)

# 2) Update the first paragraph's formatting: add a paragraph border
#    (space-only, no visible line) and widen the left indent.
$p1 = $d.Paragraphs(1)
$pf = $p1.Range.ParagraphFormat
$pf.Borders.DistanceFromTop = 5
$pf.Borders.DistanceFromLeft = 5
$pf.Borders.DistanceFromBottom = 5
$pf.Borders.DistanceFromRight = 5
$pf.LeftIndent = 11.25
